$d = $word.ActiveDocument
$newText = "môžete pozorovať súhvezdie Cygnus: 10. in 19. avgust, 9. in 18. septembra, 8. in 17. oktober"

# Target paragraph indices (in the original document) containing the old
# "2018: Datumi kampanje za opazovanje Perseus: ..." text that must be
# replaced wholesale by a single, unformatted run with the new text.
# Processed from last to first so earlier indices stay valid.
$targets = @(128, 86, 54, 3)

foreach ($idx in $targets) {
    $p = $d.Paragraphs($idx)
    $r = $p.Range

    # Delete the paragraph's run content but keep the trailing paragraph
    # mark (End - 1) so the paragraph itself (and its pPr) survives.
    $delRange = $d.Range($r.Start, $r.End - 1)
    $delRange.Delete()

    # Insert the new text as a single, plain run (no inherited rPr) at the
    # now-empty paragraph's start.
    $insPoint = $d.Range($p.Range.Start, $p.Range.Start)
    $insPoint.InsertAfter($newText)
}
